$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I")

# ---------------------------------------------------------------------------
# 1. Insert a new row above (old) row 4. This shifts:
#      old row 4 ("Number of disability persons" data)            -> row 5
#      old row 5 (merged footnote / Source row)                   -> row 6
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. Title (row 1) - new wording, merged across A1:I1, bold Arial 11, centered
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Martvili Municipality"
$ws.Range("A1:I1").Merge()
$titleRange = $ws.Range("A1:I1")
$titleRange.RowHeight = 51
$titleRange.HorizontalAlignment = -4108   # xlCenter
$titleRange.VerticalAlignment = -4108     # xlCenter
$titleRange.WrapText = $true
$titleRange.Font.Name = "Arial"
$titleRange.Font.Size = 11
$titleRange.Font.Bold = $true
$titleRange.Interior.Pattern = -4142      # xlPatternNone
$titleRange.Borders.LineStyle = -4142     # xlLineStyleNone

# ---------------------------------------------------------------------------
# 3. Row 2 subtitle - text unchanged, font/fill preserved, default row height
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "(End of year, persons)"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("A2").Font.Bold = $false
$ws.Range("A2").Interior.ThemeColor = 0
$ws.Range("A2").Interior.Pattern = 1
$ws.Range("A2").Borders.LineStyle = -4142
$ws.Range("A2").HorizontalAlignment = -4131  # xlLeft (general-ish, default)
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Row 3 (year header row) - A3 restyled: Arial 10 (automatic colour), no
#    fill, thin top border only. Year values (B3:I3) keep their formatting.
# ---------------------------------------------------------------------------
$a3 = $ws.Range("A3")
$a3.Font.Name = "Arial"
$a3.Font.Size = 10
$a3.Font.ColorIndex = -4105     # xlColorIndexAutomatic
$a3.Interior.Pattern = -4142    # xlPatternNone
$a3.Borders.Item(7).LineStyle = -4142   # left none
$a3.Borders.Item(8).LineStyle = 1       # top thin
$a3.Borders.Item(9).LineStyle = -4142   # bottom none
$a3.Borders.Item(10).LineStyle = -4142  # right none

# ---------------------------------------------------------------------------
# 5. Row 4 (new) - "family with disabilities Persons" label + totals
# ---------------------------------------------------------------------------
$a4 = $ws.Range("A4")
$a4.Value = "family with disabilities Persons "
$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.ThemeColor = 1
$a4.Font.Bold = $false
$a4.Interior.ThemeColor = 0
$a4.Interior.Pattern = 1
$a4.HorizontalAlignment = -4131   # xlLeft
$a4.VerticalAlignment = -4108     # xlCenter
$a4.WrapText = $true
$a4.Borders.Item(7).LineStyle = -4142
$a4.Borders.Item(8).LineStyle = 1        # top thin
$a4.Borders.Item(9).LineStyle = -4142
$a4.Borders.Item(10).LineStyle = -4142

$row4vals = @(1202,1176,1134,1175,1189,1210,1214,1234)
for ($i = 0; $i -lt 8; $i++) {
    $addr = $cols[$i] + "4"
    $cell = $ws.Range($addr)
    $cell.Value = $row4vals[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = -4105
    $cell.Interior.ThemeColor = 0
    $cell.Interior.Pattern = 1
    $cell.HorizontalAlignment = 1     # xlGeneral
    $cell.Borders.LineStyle = -4142
}
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 6. Row 5 (was row 4) - "disabilities Persons" label + refreshed totals
# ---------------------------------------------------------------------------
$a5 = $ws.Range("A5")
$a5.Value = "disabilities Persons "
$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.ThemeColor = 1
$a5.Font.Bold = $false
$a5.Interior.ThemeColor = 0
$a5.Interior.Pattern = 1
$a5.HorizontalAlignment = -4131
$a5.VerticalAlignment = -4108
$a5.WrapText = $true
$a5.Borders.Item(7).LineStyle = -4142
$a5.Borders.Item(8).LineStyle = -4142
$a5.Borders.Item(9).LineStyle = 1        # bottom thin
$a5.Borders.Item(10).LineStyle = -4142

$row5vals = @(1412,1375,1334,1374,1386,1416,1426,1457)
for ($i = 0; $i -lt 8; $i++) {
    $addr = $cols[$i] + "5"
    $cell = $ws.Range($addr)
    $cell.Value = $row5vals[$i]
    $cell.NumberFormat = "#\ ##0"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = -4105
    $cell.Interior.ThemeColor = 0
    $cell.Interior.Pattern = 1
    $cell.HorizontalAlignment = 1
    $cell.Borders.LineStyle = -4142
}
# Last cell (I5) gets a thin bottom border, matching the table's closing rule
$ws.Range("I5").Borders.Item(9).LineStyle = 1

$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------------
# 7. Row 6 (was row 5) - merged footnote / Source row, content unchanged,
#    A6 loses its border, B6:H6 keep a thin top border. New row height.
# ---------------------------------------------------------------------------
$a6 = $ws.Range("A6")
$a6.Font.Name = "Arial"
$a6.Font.Size = 9
$a6.Font.ColorIndex = -4105
$a6.Interior.ThemeColor = 0
$a6.Interior.Pattern = 1
$a6.HorizontalAlignment = -4131
$a6.VerticalAlignment = -4108
$a6.WrapText = $true
$a6.Borders.LineStyle = -4142

$restRow6 = $ws.Range("B6:H6")
$restRow6.Font.Name = "Arial"
$restRow6.Font.Size = 9
$restRow6.Font.ColorIndex = -4105
$restRow6.Interior.ThemeColor = 0
$restRow6.Interior.Pattern = 1
$restRow6.HorizontalAlignment = -4131
$restRow6.VerticalAlignment = -4108
$restRow6.WrapText = $true
$restRow6.Borders.Item(8).LineStyle = 1   # top thin

$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Column A width
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.81640625

Write-Output "done"
